$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 529.8333
$ws.Range("I19").Value = 501
$ws.Range("J19").Value = 535.6
$ws.Range("K19").Value = 501
$ws.Range("L19").Value = 535.6
$ws.Range("M19").Value = -326
$ws.Range("N19").Value = -885.6

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 142859360
$ws.Range("I64").Value = 250001800
$ws.Range("J64").Value = 2783.3333
$ws.Range("K64").Value = 250001800
$ws.Range("L64").Value = 2783.3333
$ws.Range("M64").Value = -250001552
$ws.Range("N64").Value = -3279.3333

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 142859360
$ws.Range("I67").Value = 250001800
$ws.Range("J67").Value = 2783.3333
$ws.Range("K67").Value = 250001800
$ws.Range("L67").Value = 2783.3333
$ws.Range("M67").Value = -250000942
$ws.Range("N67").Value = -4499.3333

# ALC row 82
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 600
$ws.Range("I82").Value = 600
$ws.Range("K82").Value = 1800
$ws.Range("M82").Value = -1394

# ALC row 85
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 600
$ws.Range("I85").Value = 600
$ws.Range("K85").Value = 1800
$ws.Range("M85").Value = -396

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2670.4924
$ws.Range("I132").Value = 2443.9348
$ws.Range("K132").Value = 7331.8044
$ws.Range("M132").Value = -4801.8044

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4527.4185
$ws.Range("I138").Value = 1584.92
$ws.Range("J138").Value = 8614.223
$ws.Range("K138").Value = 4754.76
$ws.Range("L138").Value = 25842.669
$ws.Range("M138").Value = 385.2399999999998
$ws.Range("N138").Value = -36122.669

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 32260316
$ws.Range("I61").Value = 40002156
$ws.Range("J61").Value = 2653.6667
$ws.Range("K61").Value = 40002156
$ws.Range("L61").Value = 2653.6667
$ws.Range("M61").Value = -40001944
$ws.Range("N61").Value = -3077.6667

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# ARM row 81
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# ARM row 84
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 32260316
$ws.Range("I136").Value = 40002156
$ws.Range("J136").Value = 2653.6667
$ws.Range("K136").Value = 120006468
$ws.Range("L136").Value = 7961.000100000001
$ws.Range("M136").Value = -120003918
$ws.Range("N136").Value = -13061.0001

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3500
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 3500
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -6496

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 576.75
$ws.Range("I107").Value = 576.75
$ws.Range("K107").Value = 576.75
$ws.Range("M107").Value = 1343.25

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2898.3125
$ws.Range("I134").Value = 2963.1
$ws.Range("K134").Value = 8889.299999999999
$ws.Range("M134").Value = -6354.299999999999

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2420
$ws.Range("I81").Value = 400
$ws.Range("J81").Value = 2925
$ws.Range("K81").Value = 1200
$ws.Range("L81").Value = 8775
$ws.Range("M81").Value = -77
$ws.Range("N81").Value = -11021

# CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 2420
$ws.Range("I84").Value = 400
$ws.Range("J84").Value = 2925
$ws.Range("K84").Value = 3600
$ws.Range("L84").Value = 26325
$ws.Range("M84").Value = 2016
$ws.Range("N84").Value = -37557

# GSM row 93
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -38744

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6350.5713
$ws.Range("I102").Value = 6125
$ws.Range("J102").Value = 6651.3335
$ws.Range("K102").Value = 6125
$ws.Range("L102").Value = 6651.3335
$ws.Range("M102").Value = -4503
$ws.Range("N102").Value = -9895.333500000001

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1403.3715
$ws.Range("I68").Value = 1248.5
$ws.Range("J68").Value = 2022.8572
$ws.Range("K68").Value = 1248.5
$ws.Range("L68").Value = 2022.8572
$ws.Range("M68").Value = -499.5
$ws.Range("N68").Value = -3520.8572

# LTW row 69
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 39800
$ws.Range("J69").Value = 39800
$ws.Range("L69").Value = 39800
$ws.Range("N69").Value = -41422

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1403.3715
$ws.Range("I71").Value = 1248.5
$ws.Range("J71").Value = 2022.8572
$ws.Range("K71").Value = 6242.5
$ws.Range("L71").Value = 10114.286
$ws.Range("M71").Value = -2498.5
$ws.Range("N71").Value = -17602.286

# LTW row 72
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H72").Value = 39800
$ws.Range("J72").Value = 39800
$ws.Range("L72").Value = 119400
$ws.Range("N72").Value = -127512

# LTW row 139
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 40715
$ws.Range("J139").Value = 40715
$ws.Range("L139").Value = 40715
$ws.Range("N139").Value = -50995

# WVR row 80
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 35966.332
$ws.Range("J80").Value = 35966.332
$ws.Range("L80").Value = 35966.332
$ws.Range("N80").Value = -37962.332

# WVR row 83
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 35966.332
$ws.Range("J83").Value = 35966.332
$ws.Range("L83").Value = 107898.996
$ws.Range("N83").Value = -117882.996

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 100500000
$ws.Range("I122").Value = 100500000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 301500000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -301497550
$ws.Range("N122").ClearContents()

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4693.4243
$ws.Range("I136").Value = 6154.4
$ws.Range("J136").Value = 2445.7693
$ws.Range("K136").Value = 18463.2
$ws.Range("L136").Value = 7337.3079
$ws.Range("M136").Value = -15913.2
$ws.Range("N136").Value = -12437.3079
